# fix 'sd' and 'mom' (mu, mom) init value
#
# Rewrites the second training-run series ("mom"/"sd" column) on Sheet1
# (H/I columns) and Sheet2 (K/M columns, loss-chart source), plus the
# final accuracy cell and the three elapsed-time cells, which are shared
# between both sheets. Also tightens the accuracy-chart's value-axis
# minimum and restores each sheet's active-cell selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# New "###.###### acc=" labels for Sheet1!H2:H31 (also rebuilds the
# shared-string table so Adam/SDProp/etc. compact down first).
$H = @(
    "607.50793 acc=",
    "25.119167 acc=",
    "25.625402 acc=",
    "19.004848 acc=",
    "25.825266 acc=",
    "15.177422 acc=",
    "13.506539 acc=",
    "8.753591 acc=",
    "8.180258 acc=",
    "4.885234 acc=",
    "9.408338 acc=",
    "8.384188 acc=",
    "13.498153 acc=",
    "5.162536 acc=",
    "4.1423407 acc=",
    "6.944345 acc=",
    "2.6742256 acc=",
    "1.323509 acc=",
    "1.9286346 acc=",
    "3.0258079 acc=",
    "2.185761 acc=",
    "2.2570162 acc=",
    "2.565874 acc=",
    "3.8926957 acc=",
    "3.4811668 acc=",
    "8.755626 acc=",
    "2.6696956 acc=",
    "8.219879 acc=",
    "3.756762 acc=",
    "2.5063317 acc="
)

# New accuracy values, Sheet1!I2:I31 (== Sheet2!M2:M31).
$I = @(
    0.10100000000000001,
    0.89770000000000005,
    0.92210000000000003,
    0.93720000000000003,
    0.94220000000000004,
    0.95199999999999996,
    0.95640000000000003,
    0.96030000000000004,
    0.96,
    0.9657,
    0.9698,
    0.97150000000000003,
    0.97270000000000001,
    0.97019999999999995,
    0.97430000000000005,
    0.97450000000000003,
    0.97829999999999995,
    0.97819999999999996,
    0.98019999999999996,
    0.97829999999999995,
    0.9798,
    0.97940000000000005,
    0.98229999999999995,
    0.98089999999999999,
    0.9819,
    0.98080000000000001,
    0.98460000000000003,
    0.98380000000000001,
    0.98409999999999997,
    0.98240000000000005
)

# New loss values, Sheet2!K2:K31 (full precision; H-string above is the
# rounded display copy of the same number).
$K = @(
    607.50792999999999,
    25.119167000000001,
    25.625402000000001,
    19.004847999999999,
    25.825265999999999,
    15.177422,
    13.506539,
    8.7535910000000001,
    8.1802580000000003,
    4.8852339999999996,
    9.4083380000000005,
    8.384188,
    13.498153,
    5.1625360000000002,
    4.1423407000000001,
    6.9443450000000002,
    2.6742256000000002,
    1.323509,
    1.9286346000000001,
    3.0258079000000002,
    2.1857609999999998,
    2.2570161999999998,
    2.565874,
    3.8926957,
    3.4811668,
    8.7556259999999995,
    2.6696955999999998,
    8.2198790000000006,
    3.7567620000000002,
    2.5063317000000001
)

# 1) Sheet1 H2:H31 / I2:I31  (write strings first so the shared-string
#    table compacts/reorders the same way as the reference edit).
for ($i = 0; $i -lt $H.Length; $i++) {
    $row = $i + 2
    $ws1.Range("H$row").Value = $H[$i]
    $ws1.Range("I$row").Value = $I[$i]
}

# 2) Sheet2 K2:K31 / M2:M31
for ($i = 0; $i -lt $K.Length; $i++) {
    $row = $i + 2
    $ws2.Range("K$row").Value = $K[$i]
    $ws2.Range("M$row").Value = $I[$i]
}

# 3) Final accuracy cell, shared value on both sheets.
$ws1.Range("G32").Value = 0.98509999999999998
$ws2.Range("I32").Value = 0.98509999999999998

# 4) Elapsed-time strings, shared text between Sheet1!G34:G36 and
#    Sheet2!I34:I36.
$ws1.Range("G34").Value = "5m35.011s"
$ws1.Range("G35").Value = "15m51.091s"
$ws1.Range("G36").Value = "2m34.909s"

$ws2.Range("I34").Value = "5m35.011s"
$ws2.Range("I35").Value = "15m51.091s"
$ws2.Range("I36").Value = "2m34.909s"

# 5) Accuracy chart (Sheet1) value axis: floor the scale at 0.8.
$chartObj = $ws1.ChartObjects(1)
$chart = $chartObj.Chart
$valAx = $chart.Axes(2)
$valAx.MinimumScale = 0.8

# 6) Restore per-sheet selections, leaving Sheet1 as the active tab.
$ws2.Activate()
$ws2.Range("P34").Select()

$ws1.Activate()
$ws1.Range("F2").Select()
